# Update the team-specific game outcome matrix on Sheet1 (Apprentice_B) with
# refreshed probabilities after simulating more games.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2592592592592592
$ws.Range("C2").Value = 0.4074074074074074
$ws.Range("J2").Value = 0.1111111111111111
$ws.Range("P2").Value = 0.1851851851851852
$ws.Range("S2").Value = 0.03703703703703703
$ws.Range("J3").Value = 0.09090909090909091
$ws.Range("P3").Value = 0.5454545454545454
$ws.Range("S3").Value = 0.3636363636363636
$ws.Range("J4").Value = 0.25
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.25
$ws.Range("J5").Value = 1
$ws.Range("D6").Value = 0.05882352941176471
$ws.Range("J6").Value = 0.1764705882352941
$ws.Range("O6").Value = 0.05882352941176471
$ws.Range("Q6").Value = 0.1176470588235294
$ws.Range("R6").Value = 0.1764705882352941
$ws.Range("S6").Value = 0.4117647058823529
$ws.Range("J7").Value = 0.1666666666666667
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.3333333333333333
$ws.Range("S7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.1176470588235294
$ws.Range("F8").Value = 0.05882352941176471
$ws.Range("J8").Value = 0.4117647058823529
$ws.Range("Q8").Value = 0.1764705882352941
$ws.Range("S8").Value = 0.2352941176470588
$ws.Range("B9").Value = 0.05882352941176471
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.05882352941176471
$ws.Range("Q9").Value = 0.2941176470588235
$ws.Range("R9").Value = 0.2352941176470588
$ws.Range("S9").Value = 0.2941176470588235
$ws.Range("B10").Value = 0.1134751773049645
$ws.Range("D10").Value = 0.01418439716312057
$ws.Range("F10").Value = 0.07092198581560284
$ws.Range("J10").Value = 0.148936170212766
$ws.Range("O10").Value = 0.01418439716312057
$ws.Range("Q10").Value = 0.2836879432624114
$ws.Range("R10").Value = 0.04964539007092199
$ws.Range("S10").Value = 0.3049645390070922
$ws.Range("G11").Value = 0.1578947368421053
$ws.Range("J11").Value = 0.2105263157894737
$ws.Range("K11").Value = 0.3157894736842105
$ws.Range("L11").Value = 0.3157894736842105
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.3333333333333333
$ws.Range("H15").Value = 0.1176470588235294
$ws.Range("J15").Value = 0.4705882352941176
$ws.Range("O15").Value = 0.2941176470588235
$ws.Range("S15").Value = 0.1176470588235294
$ws.Range("H16").Value = 0.1538461538461539
$ws.Range("J16").Value = 0.8461538461538461
$ws.Range("F17").Value = 0.02040816326530612
$ws.Range("H17").Value = 0.1020408163265306
$ws.Range("I17").Value = 0.1428571428571428
$ws.Range("J17").Value = 0.5102040816326531
$ws.Range("K17").Value = 0.08163265306122448
$ws.Range("O17").Value = 0.1020408163265306
$ws.Range("S17").Value = 0.04081632653061224
$ws.Range("H18").Value = 0.07142857142857142
$ws.Range("I18").Value = 0.2142857142857143
$ws.Range("J18").Value = 0.5714285714285714
$ws.Range("O18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01351351351351351
$ws.Range("H19").Value = 0.0945945945945946
$ws.Range("I19").Value = 0.08108108108108109
$ws.Range("J19").Value = 0.6081081081081081
$ws.Range("K19").Value = 0.1216216216216216
$ws.Range("O19").Value = 0.01351351351351351
$ws.Range("S19").Value = 0.06756756756756757
